# "Generate Report for Handback" - mark a.md as handed-back (in sync with
# en-US) for both locales, and record the latest target/handback file +
# datetime for the zh-cn and de-de handback rows.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# --- Overview sheet: roll the per-locale status text up for a.md (row 2)
# and b.md (row 3) ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# --- zh-cn sheet: Status column + newly generated target/handback info for
# a.md (row 2) and b.md (row 3) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusText
$wsZh.Range("C3").Value = $statusText

$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K2").Value = "2016-08-16 20:33:40"

$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("K3").Value = "2016-08-16 20:33:40"

# rebuild hyperlinks: a.md (A2), the new target link (I2), b.md (A3), the
# new target link (I3) - same target URLs as the existing a.md / b.md links
$zhA2Url = $wsZh.Hyperlinks.Item(1).Address
$zhA3Url = $wsZh.Hyperlinks.Item(2).Address
$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhA2Url, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $zhA2Url, "", "", "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $zhA3Url, "", "", "b.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $zhA2Url, "", "", "a.md")

# --- de-de sheet: same shape as zh-cn, different xlf name + handback time ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusText
$wsDe.Range("C3").Value = $statusText

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K2").Value = "2016-08-16 20:33:47"

$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("K3").Value = "2016-08-16 20:33:47"

$deA2Url = $wsDe.Hyperlinks.Item(1).Address
$deA3Url = $wsDe.Hyperlinks.Item(2).Address
$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deA2Url, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $deA2Url, "", "", "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $deA3Url, "", "", "b.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $deA2Url, "", "", "a.md")

# --- Column widths: the longer status text / new target-file column need
# more room (values picked so the engine's column-width quantization lands
# on the same stored width the original report shows) ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.16666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.16666667

$wsZh.Columns.Item(3).ColumnWidth = 29.16666667
$wsZh.Columns.Item(10).ColumnWidth = 39.16666667

$wsDe.Columns.Item(3).ColumnWidth = 29.16666667
$wsDe.Columns.Item(10).ColumnWidth = 39.16666667
